# Auto-generated Excel COM-interop edit script
# Applies the numeric updates to columns H-N across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 12 (Leve Item ID=5515)
$ws.Range("H12").Value = 1148
$ws.Range("I12").Value = 1283.4286
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 1283.4286
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -1113.4286
$ws.Range("N12").Value = -540
# row 21 (Leve Item ID=2149)
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35936
# row 23 (Leve Item ID=2149)
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 35000
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35468
# row 86 (Leve Item ID=12603)
$ws.Range("H86").Value = 9118.615
$ws.Range("I86").Value = 1388.8889
$ws.Range("J86").Value = 26510.5
$ws.Range("K86").Value = 1388.8889
$ws.Range("L86").Value = 26510.5
$ws.Range("M86").Value = -265.8888999999999
$ws.Range("N86").Value = -28756.5
# row 89 (Leve Item ID=12603)
$ws.Range("H89").Value = 9118.615
$ws.Range("I89").Value = 1388.8889
$ws.Range("J89").Value = 26510.5
$ws.Range("K89").Value = 6944.4445
$ws.Range("L89").Value = 132552.5
$ws.Range("M89").Value = -1328.4445
$ws.Range("N89").Value = -143784.5
# row 129 (Leve Item ID=36115)
$ws.Range("H129").Value = 324433.3
$ws.Range("J129").Value = 359165.56
$ws.Range("L129").Value = 1077496.68
$ws.Range("N129").Value = -1087496.68
# row 137 (Leve Item ID=44013)
$ws.Range("H137").Value = 72715.44500000001
$ws.Range("I137").Value = 103733.586
$ws.Range("J137").Value = 1556.1765
$ws.Range("K137").Value = 311200.758
$ws.Range("L137").Value = 4668.529500000001
$ws.Range("M137").Value = -308650.758
$ws.Range("N137").Value = -9768.529500000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID=27713)
$ws.Range("H2").Value = 1060.1852
$ws.Range("I2").Value = 1017.3158
$ws.Range("J2").Value = 1162
$ws.Range("K2").Value = 1017.3158
$ws.Range("L2").Value = 1162
$ws.Range("M2").Value = -904.3158
$ws.Range("N2").Value = -1388
# row 23 (Leve Item ID=2236)
$ws.Range("H23").Value = 12000
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 4000
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = -19741
$ws.Range("N23").Value = -4518
# row 32 (Leve Item ID=44147)
$ws.Range("H32").Value = 8972.057000000001
$ws.Range("I32").Value = 6563.12
$ws.Range("J32").Value = 21877.072
$ws.Range("K32").Value = 6563.12
$ws.Range("L32").Value = 21877.072
$ws.Range("M32").Value = -6276.12
$ws.Range("N32").Value = -22451.072
# row 45 (Leve Item ID=27714)
$ws.Range("H45").Value = 3482.8096
$ws.Range("I45").Value = 3209.25
$ws.Range("J45").Value = 3847.5557
$ws.Range("K45").Value = 3209.25
$ws.Range("L45").Value = 3847.5557
$ws.Range("M45").Value = -2832.25
$ws.Range("N45").Value = -4601.5557
# row 61 (Leve Item ID=43999)
$ws.Range("H61").Value = 8336154.5
$ws.Range("J61").Value = 4047.5
$ws.Range("L61").Value = 4047.5
$ws.Range("N61").Value = -4471.5
# row 74 (Leve Item ID=44000)
$ws.Range("H74").Value = 29413314
$ws.Range("I74").Value = 43478944
$ws.Range("J74").Value = 3363.6365
$ws.Range("K74").Value = 43478944
$ws.Range("L74").Value = 3363.6365
$ws.Range("M74").Value = -43478070
$ws.Range("N74").Value = -5111.636500000001
# row 77 (Leve Item ID=44000)
$ws.Range("H77").Value = 29413314
$ws.Range("I77").Value = 43478944
$ws.Range("J77").Value = 3363.6365
$ws.Range("K77").Value = 217394720
$ws.Range("L77").Value = 16818.1825
$ws.Range("M77").Value = -217390352
$ws.Range("N77").Value = -25554.1825
# row 116 (Leve Item ID=27713)
$ws.Range("H116").Value = 1060.1852
$ws.Range("I116").Value = 1017.3158
$ws.Range("J116").Value = 1162
$ws.Range("K116").Value = 1017.3158
$ws.Range("L116").Value = 1162
$ws.Range("M116").Value = 1276.6842
$ws.Range("N116").Value = -5750
# row 136 (Leve Item ID=43999)
$ws.Range("H136").Value = 8336154.5
$ws.Range("J136").Value = 4047.5
$ws.Range("L136").Value = 12142.5
$ws.Range("N136").Value = -17242.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID=27713)
$ws.Range("H3").Value = 1060.1852
$ws.Range("I3").Value = 1017.3158
$ws.Range("J3").Value = 1162
$ws.Range("K3").Value = 1017.3158
$ws.Range("L3").Value = 1162
$ws.Range("M3").Value = -903.3158
$ws.Range("N3").Value = -1390
# row 107 (Leve Item ID=27706)
$ws.Range("H107").Value = 2417.2903
$ws.Range("I107").Value = 2017.0869
$ws.Range("J107").Value = 3567.875
$ws.Range("K107").Value = 2017.0869
$ws.Range("L107").Value = 3567.875
$ws.Range("M107").Value = -97.08690000000001
$ws.Range("N107").Value = -7407.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID=44023)
$ws.Range("H31").Value = 4443.1177
$ws.Range("I31").Value = 2446.5
$ws.Range("J31").Value = 5957.793
$ws.Range("K31").Value = 2446.5
$ws.Range("L31").Value = 5957.793
$ws.Range("M31").Value = -2151.5
$ws.Range("N31").Value = -6547.793
# row 34 (Leve Item ID=44023)
$ws.Range("H34").Value = 4443.1177
$ws.Range("I34").Value = 2446.5
$ws.Range("J34").Value = 5957.793
$ws.Range("K34").Value = 2446.5
$ws.Range("L34").Value = 5957.793
$ws.Range("M34").Value = -2244.5
$ws.Range("N34").Value = -6361.793
# row 58 (Leve Item ID=44021)
$ws.Range("H58").Value = 14864.486
$ws.Range("I58").Value = 1298.6818
$ws.Range("J58").Value = 34761
$ws.Range("K58").Value = 1298.6818
$ws.Range("L58").Value = 34761
$ws.Range("M58").Value = -1095.6818
$ws.Range("N58").Value = -35167
# row 94 (Leve Item ID=32934)
$ws.Range("H94").Value = 4117.769
$ws.Range("I94").Value = 2917.1667
$ws.Range("J94").Value = 5146.857
$ws.Range("K94").Value = 2917.1667
$ws.Range("L94").Value = 5146.857
$ws.Range("M94").Value = -2466.1667
$ws.Range("N94").Value = -6048.857
# row 97 (Leve Item ID=19730)
$ws.Range("H97").Value = 32000
$ws.Range("J97").Value = 32000
$ws.Range("L97").Value = 32000
$ws.Range("N97").Value = -33982
# row 105 (Leve Item ID=19928)
$ws.Range("H105").Value = 3602.2
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253
# row 132 (Leve Item ID=44019)
$ws.Range("H132").Value = 32260182
$ws.Range("I132").Value = 35715660
$ws.Range("J132").Value = 9071
$ws.Range("K132").Value = 107146980
$ws.Range("L132").Value = 27213
$ws.Range("M132").Value = -107144450
$ws.Range("N132").Value = -32273
# row 134 (Leve Item ID=44020)
$ws.Range("H134").Value = 52632496
$ws.Range("I134").Value = 62500932
$ws.Range("J134").Value = 842.6667
$ws.Range("K134").Value = 187502796
$ws.Range("L134").Value = 2528.0001
$ws.Range("M134").Value = -187500261
$ws.Range("N134").Value = -7598.0001
# row 136 (Leve Item ID=44021)
$ws.Range("H136").Value = 14864.486
$ws.Range("I136").Value = 1298.6818
$ws.Range("J136").Value = 34761
$ws.Range("K136").Value = 3896.0454
$ws.Range("L136").Value = 104283
$ws.Range("M136").Value = -1346.0454
$ws.Range("N136").Value = -109383

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 36 (Leve Item ID=4732)
$ws.Range("H36").Value = 1601.8
$ws.Range("I36").Value = 1
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 166
# row 56 (Leve Item ID=10146)
$ws.Range("H56").Value = 6698.4614
$ws.Range("I56").Value = 6698.4614
$ws.Range("K56").Value = 6698.4614
$ws.Range("M56").Value = -6168.4614
# row 60 (Leve Item ID=4750)
$ws.Range("H60").Value = 207.14285
$ws.Range("I60").Value = 100
$ws.Range("J60").Value = 287.5
$ws.Range("K60").Value = 300
$ws.Range("L60").Value = 862.5
$ws.Range("M60").Value = -49
$ws.Range("N60").Value = -1364.5
# row 113 (Leve Item ID=27843)
$ws.Range("H113").Value = 710.5769
$ws.Range("I113").Value = 611.6
$ws.Range("J113").Value = 772.4375
$ws.Range("K113").Value = 1834.8
$ws.Range("L113").Value = 2317.3125
$ws.Range("M113").Value = 335.1999999999998
$ws.Range("N113").Value = -6657.3125
# row 122 (Leve Item ID=36078)
$ws.Range("H122").Value = 1033.1666
$ws.Range("J122").Value = 1033.1666
$ws.Range("L122").Value = 9298.499400000001
$ws.Range("N122").Value = -14198.4994
# row 131 (Leve Item ID=36060)
$ws.Range("H131").Value = 323414.84
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 323414.84
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 970244.52
$ws.Range("M131").ClearContents() | Out-Null
$ws.Range("N131").Value = -980324.52
# row 134 (Leve Item ID=44074)
$ws.Range("H134").Value = 2921.543
$ws.Range("I134").Value = 1929.6666
$ws.Range("J134").Value = 6269.125
$ws.Range("K134").Value = 5788.9998
$ws.Range("L134").Value = 18807.375
$ws.Range("M134").Value = -718.9997999999996
$ws.Range("N134").Value = -28947.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 15 (Leve Item ID=12018)
$ws.Range("H15").Value = 17500
$ws.Range("J15").Value = 17500
$ws.Range("L15").Value = 17500
$ws.Range("N15").Value = -18076
# row 70 (Leve Item ID=14146)
$ws.Range("H70").Value = 3680882
$ws.Range("I70").Value = 4599.778
$ws.Range("J70").Value = 7816700
$ws.Range("K70").Value = 4599.778
$ws.Range("L70").Value = 7816700
$ws.Range("M70").Value = -4329.778
$ws.Range("N70").Value = -7817240
# row 73 (Leve Item ID=14146)
$ws.Range("H73").Value = 3680882
$ws.Range("I73").Value = 4599.778
$ws.Range("J73").Value = 7816700
$ws.Range("K73").Value = 4599.778
$ws.Range("L73").Value = 7816700
$ws.Range("M73").Value = -3663.778
$ws.Range("N73").Value = -7818572
# row 81 (Leve Item ID=12018)
$ws.Range("H81").Value = 17500
$ws.Range("J81").Value = 17500
$ws.Range("L81").Value = 17500
$ws.Range("N81").Value = -19496
# row 84 (Leve Item ID=12018)
$ws.Range("H84").Value = 17500
$ws.Range("J84").Value = 17500
$ws.Range("L84").Value = 52500
$ws.Range("N84").Value = -62484
# row 94 (Leve Item ID=19511)
$ws.Range("H94").Value = 26448
$ws.Range("J94").Value = 26448
$ws.Range("L94").Value = 26448
$ws.Range("N94").Value = -27800
# row 97 (Leve Item ID=19940)
$ws.Range("H97").Value = 1016.4
$ws.Range("I97").Value = 364.85715
$ws.Range("J97").Value = 2536.6667
$ws.Range("K97").Value = 364.85715
$ws.Range("L97").Value = 2536.6667
$ws.Range("M97").Value = 131.14285
$ws.Range("N97").Value = -3528.6667
# row 123 (Leve Item ID=34150)
$ws.Range("H123").Value = 5673.4375
$ws.Range("I123").Value = 3093.6843
$ws.Range("J123").Value = 9443.846
$ws.Range("K123").Value = 3093.6843
$ws.Range("L123").Value = 9443.846
$ws.Range("M123").Value = -643.6842999999999
$ws.Range("N123").Value = -14343.846
# row 132 (Leve Item ID=44008)
$ws.Range("H132").Value = 4254308
$ws.Range("I132").Value = 6052011.5
$ws.Range("J132").Value = 59666.332
$ws.Range("K132").Value = 18156034.5
$ws.Range("L132").Value = 178998.996
$ws.Range("M132").Value = -18153504.5
$ws.Range("N132").Value = -184058.996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 93 (Leve Item ID=19993)
$ws.Range("H93").Value = 1513.8889
$ws.Range("I93").Value = 1426.4706
$ws.Range("K93").Value = 1426.4706
$ws.Range("M93").Value = -178.4706000000001
# row 136 (Leve Item ID=44060)
$ws.Range("H136").Value = 2393.0688
$ws.Range("I136").Value = 2393.0688
$ws.Range("K136").Value = 7179.2064
$ws.Range("M136").Value = -4629.2064

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 4 (Leve Item ID=2996)
$ws.Range("H4").Value = 5750.5
$ws.Range("I4").Value = 5001
$ws.Range("K4").Value = 5001
$ws.Range("M4").Value = -4888
# row 113 (Leve Item ID=27752)
$ws.Range("H113").Value = 1075.9615
$ws.Range("I113").Value = 1276.7778
$ws.Range("J113").Value = 624.125
$ws.Range("K113").Value = 3830.3334
$ws.Range("L113").Value = 1872.375
$ws.Range("M113").Value = -1660.3334
$ws.Range("N113").Value = -6212.375
# row 122 (Leve Item ID=36208)
$ws.Range("H122").Value = 1760
$ws.Range("I122").Value = 2085.7144
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6257.1432
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3807.1432
$ws.Range("N122").Value = -7900
